# Auto-generated Excel COM-interop script
# Applies the cell value updates described in the commit diff
# for Jogos_da_Semana_FlashScore_2024-10-13.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 2.75
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("S5").Value = 1.58
$ws.Range("X5").Value = 7.5
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 17
# Row 7
$ws.Range("G7").Value = 1.37
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 6.9
$ws.Range("J7").Value = 1.82
$ws.Range("L7").Value = 6.1
$ws.Range("V7").Value = 1.8
$ws.Range("W7").Value = 7.9
$ws.Range("Y7").Value = 8.5
$ws.Range("AE7").Value = 19
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 21
$ws.Range("AH7").Value = 45
$ws.Range("AJ7").Value = 150
$ws.Range("AK7").Value = 70
$ws.Range("AU7").Value = 8
$ws.Range("AV7").Value = 70
$ws.Range("AW7").Value = 8
$ws.Range("AX7").Value = 37
# Row 8
$ws.Range("G8").Value = 3.7
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 1.9
$ws.Range("J8").Value = 4.1
$ws.Range("R8").Value = 1.88
$ws.Range("X8").Value = 21
$ws.Range("AD8").Value = 6.6
$ws.Range("AF8").Value = 60
$ws.Range("AG8").Value = 7.6
$ws.Range("AJ8").Value = 16.5
$ws.Range("AK8").Value = 15
$ws.Range("AL8").Value = 25
$ws.Range("AM8").Value = 450
$ws.Range("AN8").Value = 5.5
$ws.Range("AO8").Value = 20
$ws.Range("AS8").Value = 350
$ws.Range("AU8").Value = 7
$ws.Range("AY8").Value = 18
$ws.Range("BB8").Value = 250
# Row 10
$ws.Range("G10").Value = 3.6
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 1.91
$ws.Range("J10").Value = 4.5
$ws.Range("L10").Value = 2.63
$ws.Range("N10").Value = 9
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.73
$ws.Range("X10").Value = 19
$ws.Range("AH10").Value = 8.5
$ws.Range("AK10").Value = 17
$ws.Range("AO10").Value = 23
$ws.Range("AQ10").Value = 81
$ws.Range("AV10").Value = 67
# Row 13
$ws.Range("G13").Value = 2.45
$ws.Range("I13").Value = 2.8
$ws.Range("W13").Value = 7.5
$ws.Range("Z13").Value = 23
$ws.Range("AA13").Value = 21
$ws.Range("AX13").Value = 17
# Row 18
$ws.Range("G18").Value = 2.2
$ws.Range("I18").Value = 3.9
$ws.Range("J18").Value = 3.1
$ws.Range("W18").Value = 5.5
$ws.Range("X18").Value = 9
$ws.Range("Y18").Value = 11
$ws.Range("Z18").Value = 21
$ws.Range("AG18").Value = 7.5
$ws.Range("AZ18").Value = 81
# Row 19
$ws.Range("G19").Value = 2.25
$ws.Range("H19").Value = 3.1
$ws.Range("J19").Value = 3
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 9
$ws.Range("O19").Value = 1.33
$ws.Range("P19").Value = 3.25
$ws.Range("Q19").Value = 2.08
$ws.Range("R19").Value = 1.73
$ws.Range("S19").Value = 1.44
$ws.Range("T19").Value = 2.63
$ws.Range("U19").Value = 1.8
$ws.Range("V19").Value = 1.91
$ws.Range("W19").Value = 7.5
$ws.Range("Y19").Value = 9.5
$ws.Range("AC19").Value = 8.5
$ws.Range("AE19").Value = 15
$ws.Range("AF19").Value = 51
$ws.Range("AG19").Value = 9.5
$ws.Range("AK19").Value = 29
$ws.Range("AL19").Value = 41
$ws.Range("AM19").Value = 251
$ws.Range("AO19").Value = 13
$ws.Range("AP19").Value = 23
$ws.Range("AR19").Value = 67
$ws.Range("AT19").Value = 2.63
$ws.Range("AY19").Value = 29
# Row 21
$ws.Range("G21").Value = 3.5
$ws.Range("I21").Value = 2.25
$ws.Range("J21").Value = 4.33
$ws.Range("L21").Value = 3
$ws.Range("S21").Value = 1.57
$ws.Range("T21").Value = 2.25
$ws.Range("Z21").Value = 41
$ws.Range("AT21").Value = 2.25
$ws.Range("AZ21").Value = 41
# Row 22
$ws.Range("M22").Value = 1.07
$ws.Range("N22").Value = 9
$ws.Range("Q22").Value = 2.1
$ws.Range("R22").Value = 1.7
# Row 23
$ws.Range("BD23").Value = 176
# Row 24
$ws.Range("G24").Value = 1.72
$ws.Range("I24").Value = 4.05
$ws.Range("J24").Value = 2.32
$ws.Range("L24").Value = 4.45
$ws.Range("M24").Value = 1.04
$ws.Range("N24").Value = 8.25
$ws.Range("O24").Value = 1.23
$ws.Range("P24").Value = 3.75
$ws.Range("Q24").Value = 1.7
$ws.Range("R24").Value = 2.07
$ws.Range("S24").Value = 1.37
$ws.Range("T24").Value = 2.85
$ws.Range("U24").Value = 1.7
$ws.Range("V24").Value = 2.05
$ws.Range("W24").Value = 8
$ws.Range("X24").Value = 9
$ws.Range("AA24").Value = 13
$ws.Range("AB24").Value = 22
$ws.Range("AC24").Value = 8.25
$ws.Range("AE24").Value = 14
$ws.Range("AF24").Value = 55
$ws.Range("AG24").Value = 13.5
$ws.Range("AH24").Value = 25
$ws.Range("AJ24").Value = 65
$ws.Range("AL24").Value = 37
$ws.Range("AM24").Value = 400
$ws.Range("AN24").Value = 3.7
$ws.Range("AP24").Value = 17.5
$ws.Range("AQ24").Value = 29
$ws.Range("AR24").Value = 60
$ws.Range("AS24").Value = 200
$ws.Range("AT24").Value = 2.85
$ws.Range("AU24").Value = 7.3
$ws.Range("AV24").Value = 65
$ws.Range("AW24").Value = 6
$ws.Range("AX24").Value = 23
$ws.Range("AY24").Value = 28
$ws.Range("BA24").Value = 150
$ws.Range("BB24").Value = 350
